{"js": "// Apply the resume content-update edit described by the commit:\n// \"Updated resume files, resume iframe link\"\n//\n// This performs a series of exact text replacements (summary, skills\n// lists, and a few bullet/description paragraphs) plus adds a bookmark\n// that Word's authoring tool happened to insert at the start of the\n// \"Alarm Clock -\" project heading.\n\nasync function replaceText(context, oldText, newText) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"Text not found: \" + oldText.substring(0, 60));\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n\n// 1. Professional summary paragraph.\nawait replaceText(\n  context,\n  \"A Software Developer that has more than 2 years of experience building projects with a focus on front-end web development. Has strong programming fundamentals and the ability to learn new things quickly and independently. Can plan, read, write, and review code effectively. A responsible worker who can communicate and work with a team virtually or face-to-face. Very proficient with computers. Native English speaker who can speak and write in Spanish with limited working capacity. Plays guitar and piano.\",\n  \"A Software Developer that has more than 2 years of experience building projects with a focus on front-end web development using React. Has strong programming fundamentals and the ability and drive to learn new things quickly and independently. Can plan, read, write, and review code effectively. A responsible worker who can communicate well asynchronously and work with a team virtually. Interested in growing into a full-stack development position.\"\n);\n\n// 2. \"Languages and frameworks\" skill list.\nawait replaceText(\n  context,\n  \" HTML, CSS, SASS/SCSS, JavaScript, React, C++, Python, SQL\",\n  \" HTML, CSS, SASS, JavaScript, Express, React, Python, SQL\"\n);\n\n// 3. \"Tools and IDEs\" skill list.\nawait replaceText(\n  context,\n  \": Git, GitHub, Visual Studio, Visual Studio Code, Adobe Photoshop, Postman\",\n  \": Git, GitHub, Visual Studio, Visual Studio Code, Adobe Photoshop, Postman, Figma\"\n);\n\n// 4. \"Operating Systems:\" label -> \"Concepts:\" label.\nawait replaceText(context, \"Operating Systems:\", \"Concepts:\");\n\n// 5. The operating-systems value text -> the new concepts list.\nawait replaceText(\n  context,\n  \" Microsoft Windows, Linux, Android\",\n  \" BEM, Functional Programming, Object-Oriented Programming, Agile (Scrum), Gitflow, Semantic HTML, Accessibility (A11y, ARIA)\"\n);\n\n// 6. Alarm Clock project description.\nawait replaceText(\n  context,\n  \"A GUI program written in Python with user-created alarms, a stopwatch, a digital clock, and a timer.\",\n  \"A time-keeping program written in Python with user-created alarms stored in an SQLite3 database, a stopwatch, a digital clock, and a timer.\"\n);\n\n// 7. Education blurb.\nawait replaceText(\n  context,\n  \"Graduated with highest honors (summa cum laude) and a cumulative GPA of 4.5. Succeeded in high-level classes, turning in quality work on time. Was a member of the Cambridge AICE program, fulfilling requirements for both the standard diploma and the Cambridge AICE diploma. Earned various awards.\",\n  \"Graduated with a cumulative GPA of 4.5, earning both a standard and Cambridge AICE diploma. \"\n);\n\n// 8. Experience bullet.\nawait replaceText(\n  context,\n  \"Planned, developed, and deployed front-end applications in more than 2 programming languages. Performed code reviews for 20+ new and aspiring software developers online. Researched and fixed more than 200 software bugs. Played an integral part in all parts of the software development life cycle.\",\n  \"Planned, developed, and deployed front-end applications in more than 2 programming languages. Performed code reviews, providing feedback for 20+ new and aspiring software developers online. Diagnosed, researched, and fixed more than 200 software bugs. Played an integral part in all parts of the software development life cycle.\"\n);\n\n// 9. Bookmark that was added at the start of the \"Alarm Clock -\" project\n// heading paragraph.\nconst alarmHeading = context.document.body.search(\"Alarm Clock - \", { matchCase: true });\nalarmHeading.load(\"items\");\nawait context.sync();\nif (alarmHeading.items.length > 0) {\n  const startRange = alarmHeading.items[0].getRange(\"Start\");\n  startRange.insertBookmark(\"_vfisqgmfwmk1\");\n  await context.sync();\n}\n", "ps1": "# Apply the resume content-update edit described by the commit:\n# \"Updated resume files, resume iframe link\"\n#\n# This performs a series of exact text replacements (summary, skills\n# lists, and a few bullet/description paragraphs) plus adds a bookmark\n# that Word's authoring tool happened to insert at the start of the\n# \"Alarm Clock -\" project heading.\n\n$d = $word.ActiveDocument\n\nfunction Replace-ExactText($oldText, $newText) {\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $result = $find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2)\n    if (-not $result) {\n        throw \"Text not found: $oldText\"\n    }\n}\n\n# 1. Professional summary paragraph.\nReplace-ExactText `\n    \"A Software Developer that has more than 2 years of experience building projects with a focus on front-end web development. Has strong programming fundamentals and the ability to learn new things quickly and independently. Can plan, read, write, and review code effectively. A responsible worker who can communicate and work with a team virtually or face-to-face. Very proficient with computers. Native English speaker who can speak and write in Spanish with limited working capacity. Plays guitar and piano.\" `\n    \"A Software Developer that has more than 2 years of experience building projects with a focus on front-end web development using React. Has strong programming fundamentals and the ability and drive to learn new things quickly and independently. Can plan, read, write, and review code effectively. A responsible worker who can communicate well asynchronously and work with a team virtually. Interested in growing into a full-stack development position.\"\n\n# 2. \"Languages and frameworks\" skill list.\nReplace-ExactText `\n    \" HTML, CSS, SASS/SCSS, JavaScript, React, C++, Python, SQL\" `\n    \" HTML, CSS, SASS, JavaScript, Express, React, Python, SQL\"\n\n# 3. \"Tools and IDEs\" skill list.\nReplace-ExactText `\n    \": Git, GitHub, Visual Studio, Visual Studio Code, Adobe Photoshop, Postman\" `\n    \": Git, GitHub, Visual Studio, Visual Studio Code, Adobe Photoshop, Postman, Figma\"\n\n# 4. \"Operating Systems:\" label -> \"Concepts:\" label.\nReplace-ExactText \"Operating Systems:\" \"Concepts:\"\n\n# 5. The operating-systems value text -> the new concepts list.\nReplace-ExactText `\n    \" Microsoft Windows, Linux, Android\" `\n    \" BEM, Functional Programming, Object-Oriented Programming, Agile (Scrum), Gitflow, Semantic HTML, Accessibility (A11y, ARIA)\"\n\n# 6. Alarm Clock project description.\nReplace-ExactText `\n    \"A GUI program written in Python with user-created alarms, a stopwatch, a digital clock, and a timer.\" `\n    \"A time-keeping program written in Python with user-created alarms stored in an SQLite3 database, a stopwatch, a digital clock, and a timer.\"\n\n# 7. Education blurb.\nReplace-ExactText `\n    \"Graduated with highest honors (summa cum laude) and a cumulative GPA of 4.5. Succeeded in high-level classes, turning in quality work on time. Was a member of the Cambridge AICE program, fulfilling requirements for both the standard diploma and the Cambridge AICE diploma. Earned various awards.\" `\n    \"Graduated with a cumulative GPA of 4.5, earning both a standard and Cambridge AICE diploma. \"\n\n# 8. Experience bullet.\nReplace-ExactText `\n    \"Planned, developed, and deployed front-end applications in more than 2 programming languages. Performed code reviews for 20+ new and aspiring software developers online. Researched and fixed more than 200 software bugs. Played an integral part in all parts of the software development life cycle.\" `\n    \"Planned, developed, and deployed front-end applications in more than 2 programming languages. Performed code reviews, providing feedback for 20+ new and aspiring software developers online. Diagnosed, researched, and fixed more than 200 software bugs. Played an integral part in all parts of the software development life cycle.\"\n\n# 9. Bookmark that was added at the start of the \"Alarm Clock -\" project\n# heading paragraph.\n$bmFind = $d.Content.Find\n$bmFind.ClearFormatting()\n$bmFind.Text = \"Alarm Clock - \"\n$bmFind.Forward = $true\n$bmFind.Wrap = 1\n$bmFound = $bmFind.Execute()\nif ($bmFound) {\n    $bmRange = $bmFind.Parent.Duplicate\n    $bmRange.Collapse(1)\n    $d.Bookmarks.Add(\"_vfisqgmfwmk1\", $bmRange)\n}\n"}
